# Insert a new weekly price record for "Perejil" (Feria Lagunitas de Puerto
# Montt) just above the existing row 212, shifting the following rows down
# by one (A1:R268 -> A1:R269).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 212; rows 212..268 shift down to 213..269.
$ws.Rows.Item(212).Insert()

# Populate the newly inserted row 212 with the new observation.
$ws.Cells.Item(212, 1).Value = 4
$ws.Cells.Item(212, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(212, 3).Value = "Los Lagos"
$ws.Cells.Item(212, 4).Value = 44754
$ws.Cells.Item(212, 5).Value = 10
$ws.Cells.Item(212, 6).Value = 100112044
$ws.Cells.Item(212, 7).Value = "Perejil"
$ws.Cells.Item(212, 8).Value = "Sin especificar"
$ws.Cells.Item(212, 9).Value = "Primera"
$ws.Cells.Item(212, 10).Value = 160
$ws.Cells.Item(212, 11).Value = 5500
$ws.Cells.Item(212, 12).Value = 6000
$ws.Cells.Item(212, 13).Value = 5750
$ws.Cells.Item(212, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(212, 15).Value = "Región Metropolitana"
$ws.Cells.Item(212, 16).Value = 1917
$ws.Cells.Item(212, 17).Value = 3
$ws.Cells.Item(212, 18).Value = "Hortaliza"
